{"js": "// Update the date heading and the 25 division-problem answers in the\n// table. Replacements are applied strictly by position (document order),\n// since the same problem text (\"78\u00f77=11, 1\") appears twice with two\n// different target replacements, so a global text search/replace would\n// be ambiguous.\n\n// 1) Update the date/weekday heading paragraph (first paragraph in body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst heading = paragraphs.items[0];\nheading.load(\"text\");\nawait context.sync();\nif (heading.text.trim() === \"2023-10-26 Thursday\") {\n  heading.insertText(\"2023-10-27 Friday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the answer table cells. The table has 20 rows x 5 columns,\n// but only every 4th row (0, 4, 8, 12, 16) actually holds an answer -\n// the rows in between are blank spacer rows. Walk the grid in row-major\n// order and apply the next replacement from the list whenever a\n// non-empty cell is encountered.\nconst replacements = [\n  \"87\u00f78=10, 7\",\n  \"63\u00f78=7, 7\",\n  \"16\u00f77=2, 2\",\n  \"97\u00f76=16, 1\",\n  \"51\u00f73=17, 0\",\n  \"64\u00f72=32, 0\",\n  \"32\u00f75=6, 2\",\n  \"97\u00f79=10, 7\",\n  \"45\u00f72=22, 1\",\n  \"46\u00f75=9, 1\",\n  \"43\u00f74=10, 3\",\n  \"90\u00f78=11, 2\",\n  \"16\u00f76=2, 4\",\n  \"19\u00f79=2, 1\",\n  \"79\u00f72=39, 1\",\n  \"89\u00f74=22, 1\",\n  \"49\u00f76=8, 1\",\n  \"22\u00f78=2, 6\",\n  \"11\u00f76=1, 5\",\n  \"24\u00f77=3, 3\",\n  \"65\u00f77=9, 2\",\n  \"53\u00f79=5, 8\",\n  \"49\u00f78=6, 1\",\n  \"13\u00f78=1, 5\",\n  \"73\u00f72=36, 1\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values,rowCount\");\nawait context.sync();\n\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  const row = table.values[r];\n  for (let c = 0; c < row.length; c++) {\n    const current = row[c];\n    if (current !== \"\" && idx < replacements.length) {\n      const cell = table.getCell(r, c);\n      cell.value = replacements[idx];\n      idx++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division-problem answers in the\n# table. Replacements are applied strictly by position (document order),\n# since the same problem text (\"78\u00f77=11, 1\") appears twice with two\n# different target replacements, so a global text search/replace would\n# be ambiguous.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday heading paragraph (first paragraph in body).\n$heading = $d.Paragraphs.Item(1)\nif ($heading.Range.Text.TrimEnd(\"`r\") -eq \"2023-10-26 Thursday\") {\n    $heading.Range.Text = \"2023-10-27 Friday\"\n}\n\n# 2) Update the answer table cells. The table has 20 rows x 5 columns,\n# but only every 4th row (1, 5, 9, 13, 17 in 1-based indexing) actually\n# holds an answer - the rows in between are blank spacer rows. Walk the\n# grid in row-major order and apply the next replacement from the list\n# whenever a non-empty cell is encountered.\n$replacements = @(\n    \"87\u00f78=10, 7\",\n    \"63\u00f78=7, 7\",\n    \"16\u00f77=2, 2\",\n    \"97\u00f76=16, 1\",\n    \"51\u00f73=17, 0\",\n    \"64\u00f72=32, 0\",\n    \"32\u00f75=6, 2\",\n    \"97\u00f79=10, 7\",\n    \"45\u00f72=22, 1\",\n    \"46\u00f75=9, 1\",\n    \"43\u00f74=10, 3\",\n    \"90\u00f78=11, 2\",\n    \"16\u00f76=2, 4\",\n    \"19\u00f79=2, 1\",\n    \"79\u00f72=39, 1\",\n    \"89\u00f74=22, 1\",\n    \"49\u00f76=8, 1\",\n    \"22\u00f78=2, 6\",\n    \"11\u00f76=1, 5\",\n    \"24\u00f77=3, 3\",\n    \"65\u00f77=9, 2\",\n    \"53\u00f79=5, 8\",\n    \"49\u00f78=6, 1\",\n    \"13\u00f78=1, 5\",\n    \"73\u00f72=36, 1\"\n)\n\n$table = $d.Tables.Item(1)\n$idx = 0\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cellRange = $cell.Range\n        $text = $cellRange.Text.TrimEnd([char]7).TrimEnd(\"`r\")\n        if ($text -ne \"\" -and $idx -lt $replacements.Length) {\n            $cellRange.Text = $replacements[$idx]\n            $idx++\n        }\n    }\n}\n"}
